$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CtByBldg")

$ws.Range("I16").Value = 1
$ws.Range("I45").Value = 1
$ws.Range("I62").Value = 1
$ws.Range("H65").Value = 35
$ws.Range("I65").Value = 14
$ws.Range("H66").Value = 12
$ws.Range("I66").Value = 13
$ws.Range("H67").Value = 22
$ws.Range("I67").Value = 10
$ws.Range("H81").Value = 23
$ws.Range("I81").Value = 6
$ws.Range("H82").Value = 43
$ws.Range("I82").Value = 5
$ws.Range("H86").Value = 32
$ws.Range("I86").Value = 9
$ws.Range("I115").Value = 1
$ws.Range("I148").Value = 2
$ws.Range("I163").Value = 1
$ws.Range("I179").Value = 2
$ws.Range("I180").Value = 1
$ws.Range("I209").Value = 1
$ws.Range("I227").Value = 1
$ws.Range("H339").Value = 66
$ws.Range("I339").Value = 16
$ws.Range("I416").Value = 1
$ws.Range("I444").Value = 1
$ws.Range("H453").Value = 13
$ws.Range("I453").Value = 3
$ws.Range("H521").Value = 25
$ws.Range("I521").Value = 3
$ws.Range("I522").Value = 1
$ws.Range("H556").Value = 14
$ws.Range("I556").Value = 7
$ws.Range("H567").Value = 23
$ws.Range("I567").Value = 6
$ws.Range("H577").Value = 50
$ws.Range("I577").Value = 24
$ws.Range("H578").Value = 124
$ws.Range("I578").Value = 22
$ws.Range("H579").Value = 89
$ws.Range("I579").Value = 10
$ws.Range("H582").Value = 110
$ws.Range("I582").Value = 15
$ws.Range("H590").Value = 35
$ws.Range("I590").Value = 6
$ws.Range("H591").Value = 32
$ws.Range("I591").Value = 10
$ws.Range("H596").Value = 11
$ws.Range("I596").Value = 2
$ws.Range("H639").Value = 108
$ws.Range("I639").Value = 24
$ws.Range("H640").Value = 35
$ws.Range("I640").Value = 9
$ws.Range("H641").Value = 133
$ws.Range("I641").Value = 45
$ws.Range("H677").Value = 1
$ws.Range("I677").Value = 4
$ws.Range("I681").Value = 20
$ws.Range("H704").Value = 12
$ws.Range("I704").Value = 2
$ws.Range("I712").Value = 1
$ws.Range("H731").Value = 32
$ws.Range("I731").Value = 3
$ws.Range("H732").Value = 32
$ws.Range("I732").Value = 3
$ws.Range("H733").Value = 32
$ws.Range("I733").Value = 3
$ws.Range("H734").Value = 32
$ws.Range("I734").Value = 3
$ws.Range("G735").Value = 1
$ws.Range("H735").Value = 14
$ws.Range("I735").Value = 3
$ws.Range("H736").Value = 14
$ws.Range("I736").Value = 3
$ws.Range("H737").Value = 14
$ws.Range("H739").Value = 14
$ws.Range("I739").Value = 4
$ws.Range("H751").Value = 4
$ws.Range("I751").Value = 4
$ws.Range("I752").Value = 4
$ws.Range("H753").Value = 10
$ws.Range("I753").Value = 4
$ws.Range("H754").Value = 8
$ws.Range("I754").Value = 4
$ws.Range("H762").Value = 15
$ws.Range("I773").Value = 8
$ws.Range("I775").Value = 1
$ws.Range("I777").Value = 2
$ws.Range("H799").Value = 57
$ws.Range("I799").Value = 19
$ws.Range("H800").Value = 13
$ws.Range("I800").Value = 4
$ws.Range("H803").Value = 10
$ws.Range("I803").Value = 1
$ws.Range("H804").Value = 34
$ws.Range("I804").Value = 10
$ws.Range("H855").Value = 4
$ws.Range("I855").Value = 4
$ws.Range("H861").Value = 48
$ws.Range("I861").Value = 9
$ws.Range("H880").Value = 104
$ws.Range("I880").Value = 16
$ws.Range("H905").Value = 4
$ws.Range("I905").Value = 1
$ws.Range("H909").Value = 15
$ws.Range("I909").Value = 2
$ws.Range("H910").Value = 88
$ws.Range("I910").Value = 18
$ws.Range("H911").Value = 137
$ws.Range("I911").Value = 22
$ws.Range("I918").Value = 5
$ws.Range("H927").Value = 103
$ws.Range("I927").Value = 20
$ws.Range("H931").Value = 10
$ws.Range("I931").Value = 1
$ws.Range("H943").Value = 37
$ws.Range("I943").Value = 7
$ws.Range("H965").Value = 13
$ws.Range("I965").Value = 1
$ws.Range("H979").Value = 21
$ws.Range("I979").Value = 4
$ws.Range("H980").Value = 25
$ws.Range("I980").Value = 5
$ws.Range("I981").Value = 1
$ws.Range("H989").Value = 65
$ws.Range("I989").Value = 13
$ws.Range("H996").Value = 85
$ws.Range("I996").Value = 17
$ws.Range("H998").Value = 40
$ws.Range("I998").Value = 8
$ws.Range("I1095").Value = 1
$ws.Range("H1098").Value = 130
$ws.Range("I1098").Value = 24
$ws.Range("H1102").Value = 136
$ws.Range("I1102").Value = 29
$ws.Range("H1112").Value = 58
$ws.Range("I1112").Value = 14
$ws.Range("H1124").Value = 16
$ws.Range("H1125").Value = 4
$ws.Range("H1126").Value = 8
$ws.Range("I1126").Value = 2
$ws.Range("H1127").Value = 2
$ws.Range("I1127").Value = 1
$ws.Range("H1128").Value = 11
$ws.Range("I1128").Value = 2
$ws.Range("H1135").Value = 38
$ws.Range("I1135").Value = 9
$ws.Range("H1146").Value = 1
$ws.Range("I1146").Value = 1
$ws.Range("H1148").Value = 2
$ws.Range("I1148").Value = 2
$ws.Range("G1163").Value = 1
$ws.Range("H1163").Value = 5
$ws.Range("I1163").Value = 1
$ws.Range("H1211").Value = 86
$ws.Range("I1211").Value = 31

$ws.Range("K2").Select()
